$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap "Periodo Mora" (E) and "Valor Mora" (F) between row 16 and row 17,
# and update "Salario Basico" (G) on both rows to the new value.
$ws.Range("E16").Value = "2310"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 877803

$ws.Range("E17").Value = "2309"
$ws.Range("F17").Value = 58667
$ws.Range("G17").Value = 877803
